$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'49.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.63%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.310"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.49%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08060"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.74%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.589"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.70%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.405"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'34.17%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.71%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1324"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.81%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.75%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09474"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.39%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04655"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.71%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1043"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001315"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.29%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04196"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.47%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005899"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.28%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.350"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.21%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.464"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.62%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3505"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'4.67%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.084"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-6.59%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1367"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.37%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3089"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.86%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001311"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.99%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004264"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.62%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.11%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003520"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02732"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'9.99%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06493"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'22.49%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01047"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'76.31%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008195"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'7.02%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.72%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007763"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.26%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008636"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'14.45%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3490"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'16.16%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006612"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.74%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000746"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.43%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05590"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'37.11%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-5.31%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002089"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.43%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.43%"
$ws.Range("E51").Style = "Normal"
